$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<Label>_old" -> "<Label>_FV2410",
#    "<Label>_new" -> "<Label>_FV2504" ("diff" stays as-is).
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "K1" = "diff"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into a real Excel Table ("Table1"), keeping the
#    header row's existing look (bold / shaded / bordered / centred+wrapped)
#    exactly as it already was - stash it in a scratch area first so adding
#    the table doesn't bake a fresh header style (dxf) into the table part.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy($scratch)
$headerRange.Style = "Normal"

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.Name = "Table1"

$scratch.Copy($headerRange)
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1, top-left of the scrollable
#    pane parked at A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
